# feat: change the way to find the global best
#
# The underlying PSO run was re-executed with a new strategy for picking the
# global-best particle, which produced a fresh set of "best found" values
# (column C). This recomputes the gap (column D) / avg (column E) formulas
# automatically. A status header "Starting the run!" is also written to C1,
# and the sheet selection is left on C2:C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C.
$ws.Range("C1").Value = "Starting the run!"

# Updated "best found" results for each problem instance (column C).
$updates = [ordered]@{
    2   = 24373
    3   = 24645
    4   = 24514
    5   = 23824
    6   = 24595
    7   = 25260
    8   = 25796
    9   = 23813
    10  = 25096
    11  = 24744
    12  = 42894
    13  = 41668
    14  = 42646
    15  = 45057
    16  = 42149
    17  = 43171
    18  = 41573
    19  = 45028
    20  = 44426
    21  = 44258
    22  = 60269
    23  = 62640
    24  = 60089
    25  = 60402
    26  = 61104
    27  = 59292
    28  = 61559
    29  = 61495
    30  = 59497
    31  = 60189
    212 = 55199
    213 = 54857
    214 = 53439
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}

# Leave the sheet scrolled to the top with C2:C3 selected.
$ws.Range("C2:C3").Select()
